# Controle_Financeiro.xlsx edit
# Commit: "Cria módulo para leitura e manipulação no Excel, utilizando a biblioteca openpyxl."
#
# Summary of changes applied:
#  - ENTRADAS: header "Mês" -> "Mês " (trailing space); new empty styled (underline) cell E5
#  - GASTOS: header "Dia" -> "Dia ", "Mês" -> "Mês " (trailing spaces);
#            "Nº Parcela" -> "Numero de Parcelas"
#  - GASTOS tab becomes the active/selected sheet (selection N5); ENTRADAS keeps selection A2
#  - Column width tweaks on ENTRADAS (E:G) and GASTOS (C:J)

$wb = $excel.ActiveWorkbook

$wsEntradas = $wb.Worksheets.Item("ENTRADAS")
$wsGastos   = $wb.Worksheets.Item("GASTOS")
$wsResumo   = $wb.Worksheets.Item("RESUMO")

# ---- ENTRADAS header row ----
$wsEntradas.Range("A1").Value = "Dia"
$wsEntradas.Range("B1").Value = "Mês "
$wsEntradas.Range("C1").Value = "Ano"
$wsEntradas.Range("D1").Value = "Descrição"
$wsEntradas.Range("E1").Value = "Categoria"
$wsEntradas.Range("F1").Value = "Valor"
$wsEntradas.Range("G1").Value = "Forma de Pagamento"

# New styled (underlined) empty cell E5 -- also extends the used range to A1:G5
$wsEntradas.Range("E5").Font.Underline = $true

# ---- GASTOS header row ----
$wsGastos.Range("A1").Value = "Dia "
$wsGastos.Range("B1").Value = "Mês "
$wsGastos.Range("C1").Value = "Ano"
$wsGastos.Range("D1").Value = "Descrição"
$wsGastos.Range("E1").Value = "Categoria"
$wsGastos.Range("F1").Value = "Valor Total"
$wsGastos.Range("G1").Value = "Forma de Pagamento"
$wsGastos.Range("H1").Value = "Parcelado?"
$wsGastos.Range("I1").Value = "Numero de Parcelas"
$wsGastos.Range("J1").Value = "Total Parcelas"

# ---- RESUMO header row (values unchanged) ----
$wsResumo.Range("A1").Value = "Mês"
$wsResumo.Range("B1").Value = "Ano"
$wsResumo.Range("C1").Value = "Entradas"
$wsResumo.Range("D1").Value = "Gastos"
$wsResumo.Range("E1").Value = "Saldo"

# ---- Column widths ----
# ENTRADAS
$wsEntradas.Columns.Item(5).ColumnWidth = 17.333333333333332
$wsEntradas.Columns.Item(6).ColumnWidth = 8.666666666666666
$wsEntradas.Columns.Item(7).ColumnWidth = 18.5

# GASTOS
$wsGastos.Columns.Item(3).ColumnWidth  = 8.0
$wsGastos.Columns.Item(4).ColumnWidth  = 8.666666666666666
$wsGastos.Columns.Item(5).ColumnWidth  = 17.333333333333332
$wsGastos.Columns.Item(6).ColumnWidth  = 9.333333333333334
$wsGastos.Columns.Item(7).ColumnWidth  = 17.333333333333332
$wsGastos.Columns.Item(8).ColumnWidth  = 11.5
$wsGastos.Columns.Item(9).ColumnWidth  = 16.666666666666668
$wsGastos.Columns.Item(10).ColumnWidth = 12.5

# ---- Sheet selection / active tab: GASTOS becomes the visible/active sheet ----
$wsEntradas.Range("A2").Select() | Out-Null
$wsGastos.Activate()
$wsGastos.Range("N5").Select() | Out-Null
